# FSE Upload Template rework:
#  - rename Sheet1 -> "Format"
#  - add a new "Sample Data" sheet with header + 2 sample rows
#  - reorder/update header labels
#  - store ID-like numeric strings as Text (numFmtId 49 / "@")

$wb = $excel.ActiveWorkbook

# ---- Sheet "Format" (was "Sheet1") ----------------------------------
$format = $wb.Worksheets.Item(1)
$format.Name = "Format"

$headers = @(
    "FSE ID(Lapu no)",
    "Retailer Number",
    "Retailer Name",
    "Retailer Address",
    "Retailer Lat Long",
    "Date Of Camp (dd/mm/yyyy)",
    "Target Acquisition",
    "Target recharge count",
    "Target recharge amount",
    "Target sim Swap"
)

for ($col = 1; $col -le $headers.Length; $col++) {
    $format.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Columns B..I get re-measured ("best fit") widths for their new headers;
# A and J fall back to the sheet default (no explicit column override).
$format.Columns.Item(2).ColumnWidth = 14.0
$format.Columns.Item(3).ColumnWidth = 12.25
$format.Columns.Item(4).ColumnWidth = 13.75
$format.Columns.Item(5).ColumnWidth = 13.625
$format.Columns.Item(6).ColumnWidth = 24.25
$format.Columns.Item(7).ColumnWidth = 15.125
$format.Columns.Item(8).ColumnWidth = 19.125
$format.Columns.Item(9).ColumnWidth = 20.25

# ---- Sheet "Sample Data" (new) --------------------------------------
$sample = $wb.Worksheets.Add($null, $format)
$sample.Name = "Sample Data"

$row2 = @(
    "918377891507",
    "753011513",
    "Test Retailer",
    "Orleans",
    "7.873054/80.771797",
    "31/01/2021",
    "20",
    "20",
    "200",
    "10"
)

$row3 = @(
    "918377891507",
    "753011513",
    "Test Retailer",
    "Srilanka",
    "8.873054/82.771797",
    "30/01/2021",
    "10",
    "10",
    "200",
    "10"
)

# Mark every cell of the A1:J3 block as Text (numFmtId 49 / "@") BEFORE any
# value is written, so id-like numbers ("918377891507", "20", "200", ...)
# are stored as text rather than silently becoming numbers.
for ($row = 1; $row -le 3; $row++) {
    for ($col = 1; $col -le 10; $col++) {
        $sample.Cells.Item($row, $col).NumberFormat = "@"
    }
}

for ($col = 1; $col -le $headers.Length; $col++) {
    $sample.Cells.Item(1, $col).Value = $headers[$col - 1]
}
for ($col = 1; $col -le $row2.Length; $col++) {
    $sample.Cells.Item(2, $col).Value = $row2[$col - 1]
}
for ($col = 1; $col -le $row3.Length; $col++) {
    $sample.Cells.Item(3, $col).Value = $row3[$col - 1]
}

$sample.PageSetup.PaperSize = 9
$sample.PageSetup.Orientation = 1

# Column A is best-fit to the FSE id values, column E to the lat/long
# strings; the rest keep the sheet's standard width.
$sample.Columns.Item(1).ColumnWidth = 13.125
$sample.Columns.Item(5).ColumnWidth = 17.25

# ---- selection / active sheet -----------------------------------------
# "Sample Data" ends up with the whole of row 2 selected, "Format" (the
# tab that stays on top/active) ends up selected at D6.
$sample.Activate() | Out-Null
$sample.Range("A2:XFD2").Select() | Out-Null

$format.Activate() | Out-Null
$format.Range("D6").Select() | Out-Null
